$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("windspeed")

# Update B28:B38 values from 15 to 6
$ws.Range("B28:B38").Value = 6

# Update the selection on the sheet: sqref B2:B15, active cell B2
$ws.Activate()
$excel.Application.Goto($ws.Range("B2:B15"))
